# Commit: "Updates to NH/EF modules"
#
# The "Facets" slide (slide #4 in the deck) was removed; every slide
# after it shifts up by one position. No other content changes.

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(4)
Write-Host "Deleting slide 4:" $s.Shapes.Item(1).TextFrame.TextRange.Text
$s.Delete()
